# Apply the "Add data for 2022-07-18" update to the carjacking YoY workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (title reflects the "through" date of the latest data).
$ws.Name = "Through 2022-07-10"

# Update the label for the "July" row to reflect the new "through" date.
$ws.Range("A8").Value = "July (through 07-10)"

# Update the single-cell correction in row 5 (April 2022).
$ws.Range("I5").Value = 114

# Update row 8 (July) with new year-to-date totals.
$ws.Range("B8").Value = 13
$ws.Range("C8").Value = 15
$ws.Range("D8").Value = 17
$ws.Range("E8").Value = 27
$ws.Range("F8").Value = 14
$ws.Range("G8").Value = 34
$ws.Range("H8").Value = 53
$ws.Range("I8").Value = 53

# Update row 9 (Total) with new year-to-date totals.
$ws.Range("B9").Value = 138
$ws.Range("C9").Value = 263
$ws.Range("D9").Value = 407
$ws.Range("E9").Value = 380
$ws.Range("F9").Value = 265
$ws.Range("G9").Value = 506
$ws.Range("H9").Value = 813
$ws.Range("I9").Value = 858
